# The deck ships two embedded themes:
#   theme1.xml -> used by the (one) Slide Master  -> "Integral" / "Red Violet" palette
#   theme2.xml -> used by the Notes Master         -> "Office Theme" palette
#
# The authored change swaps the two themes' contents, so the Slide Master
# (and therefore every slide) now renders with the standard "Office" color
# palette instead of the pink/violet "Integral" one.
#
# PowerPoint's automation surface doesn't give us a "swap these two theme
# parts" verb, but it does let us repaint the twelve theme colors that back
# the deck's (single) design/theme, which is exactly the colour swap the
# diff performs on ppt/theme/theme1.xml. We do that through
# Slide.ThemeColorScheme, which edits the full modern 12-slot palette
# (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink) in one shot and - unlike the
# legacy 8-slot Master/Slide.ColorScheme - doesn't blow away the scheme's
# metadata while doing it.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$theme = $s.ThemeColorScheme

# Target palette = the "Office" theme colours (what theme1.xml becomes).
# COM RGB values are packed 0x00BBGGRR, i.e. R + G*256 + B*65536.
$theme.Item(1).RGB  = 0          # dk1      000000
$theme.Item(2).RGB  = 16777215   # lt1      FFFFFF
$theme.Item(3).RGB  = 6968388    # dk2      44546A
$theme.Item(4).RGB  = 15132391   # lt2      E7E6E6
$theme.Item(5).RGB  = 13998939   # accent1  5B9BD5
$theme.Item(6).RGB  = 3243501    # accent2  ED7D31
$theme.Item(7).RGB  = 10855845   # accent3  A5A5A5
$theme.Item(8).RGB  = 49407      # accent4  FFC000
$theme.Item(9).RGB  = 12874308   # accent5  4472C4
$theme.Item(10).RGB = 4697456    # accent6  70AD47
$theme.Item(11).RGB = 12673797   # hlink    0563C1
$theme.Item(12).RGB = 7491477    # folHlink 954F72
